$wb = $excel.ActiveWorkbook

$main = $wb.Worksheets.Item("Main")
$jesse = $wb.Worksheets.Item("Jesse")

# --- Jesse sheet: append new journal entry (row 9) ---
$jesse.Range("A9").Value = 43071
$jesse.Range("B9").Value = 360
$jesse.Range("C9").Value = "Added function getObjects() to Rooms.h to return the vector of objects. `nCreated Building.h, and Building.cpp as a Tree-like structure."
$jesse.Rows.Item(9).RowHeight = 57

# --- Update selections to match the new view state ---
$jesse.Range("H7").Select() | Out-Null
$main.Range("C18").Select() | Out-Null

# Main becomes the active/selected sheet & tab
$main.Activate()
